$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-stale conversation rows (3-8); only row 2 (the wakeword
# greeting exchange) remains below the header row.
$ws.Rows("3:8").Delete()

# Update row 2 with the new "speak on wakeword" greeting exchange.
$ws.Range("A2").Value = "2023-06-20 20:50:27"
$ws.Range("B2").Value = "hello rami"
$ws.Range("C2").Value = "hello"
$ws.Range("D2").Value = "Hello, my fantastic friend! How can I make your day shine even brighter?"
$ws.Range("E2").Value = "GEN hello"
$ws.Range("F2").Value = 2.782000000006519
$ws.Range("G2").Value = 17.5
$ws.Range("H2").Value = 20.28200000000652
